$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data to append: index (col A) and value (col B)
# (values written in plain decimal form since scientific-notation
# literals like 1E-16 are not supported by the script parser)
$newData = @(
    @(204, 0.0000000000000001295260195396016),
    @(205, 0.0000000000000001009293658750142),
    @(206, 0.00000000000000004440892098500626),
    @(207, 0.0000000000000001233581138472396),
    @(208, 0.00000000000000002775557561562891),
    @(209, 0.00000000000000003172065784643304),
    @(210, 0.0000000000000001110223024625157),
    @(211, -0.00000000000000006661338147750939),
    @(212, 0),
    @(213, -0.00000000000000005551115123125783),
    @(214, 0),
    @(215, 0)
)

$startRow = 206
$styleSource = $ws.Cells.Item($startRow - 1, 1)
for ($i = 0; $i -lt $newData.Count; $i++) {
    $row = $startRow + $i
    $pair = $newData[$i]

    # Copy the formatting (bold/border/centered style, s="1") of the
    # column-A index cells down into the newly appended rows.
    $styleSource.Copy($ws.Cells.Item($row, 1))

    $ws.Cells.Item($row, 1).Value = $pair[0]
    $ws.Cells.Item($row, 2).Value = $pair[1]
}
